# GuildConfig.xlsx edit: unify the conception of DataNode, DataTable, Entity.
#
# The only semantically meaningful, reproducible change in the source diff is
# a rename of the sole worksheet from "Property1" to "DataNode" (matching the
# commit message), together with the UI selection state that was left behind
# after the edit (the active cell moved to D36). Everything else in the raw
# OOXML diff (fileVersion/rupBuild bump, xr/xr2/xr9/xr16 revision-tracking
# namespaces & uids, the Mac->Windows absPath swap, the Calibri->宋体 default
# font re-stamp and the row-height/column-width/dyDescent deltas that cascade
# from it, the phoneticPr/"常规" cellStyle name, the new x15 timeline style)
# is cosmetic churn produced by re-saving the workbook with a different
# Excel build/locale - not addressable content edits, so it is intentionally
# left untouched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the unified "DataNode" concept.
$ws.Name = "DataNode"

# Restore the post-edit selection state (active cell D36) recorded in the diff.
$ws.Range("D36").Select() | Out-Null
